# En route to write model results
#
# 1. Bold the "EMM" label at the start of the "nota sugiere..." paragraph.
# 2. Mark the diputados.gob.mx hyperlink as visited.
# 3. Insert a new "EMM: en este caso..." paragraph (with bold "EMM") plus a
#    trailing blank paragraph after the dof.gob.mx hyperlink paragraph.
# 4. Mark the fotos.eluniversal.com.mx hyperlink as visited.

$d = $word.ActiveDocument

# --- 1) Bold "EMM" in "EMM: nota sugiere que muchos perredistas..." -------
$find1 = $d.Content
$found1 = $find1.Find.Execute("EMM: nota sugiere", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $boldPart = $d.Range($find1.Start, $find1.Start + 3)
    $boldPart.Bold = 1
}

# --- 2) diputados.gob.mx link -> VisitedInternetLink -----------------------
$hyperlinks = $d.Hyperlinks
for ($i = 1; $i -le $hyperlinks.Count; $i++) {
    $hl = $hyperlinks.Item($i)
    if ($hl.Address -eq "http://www3.diputados.gob.mx/index.php/camara/005_comunicacion/b_agencia_de_noticias/004_2008/004_abril/10_10/3745_inicio_el_fap_la_resistencia_civil_pacifica_garza_estamos_organizados_chanona_posible_sede_alterna_castano_se_cancela_la_democracia_gamboa") {
        $hl.Range.Style = "VisitedInternetLink"
    }
}

# --- 3) New "EMM: en este caso..." paragraph + trailing blank paragraph ---
$find2 = $d.Content
$found2 = $find2.Find.Execute("https://www.dof.gob.mx/nota_detalle_popup.php?codigo=5068156", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $dofPara = $find2.Paragraphs(1)
    $insertPoint = $dofPara.Range
    $insertPoint.Collapse(0)
    $insertPoint.InsertParagraphAfter()

    $newPara = $dofPara.Next()
    $newPara.Range.Text = "EMM: en este caso sí consiguieron posponer la sesión. PAN y PRI rumoraban que se reunirían en sede alterna. "

    $boldPart2 = $d.Range($newPara.Range.Start, $newPara.Range.Start + 3)
    $boldPart2.Bold = 1

    $blankInsert = $newPara.Range
    $blankInsert.Collapse(0)
    $blankInsert.InsertParagraphAfter()
}

# --- 4) fotos.eluniversal.com.mx link -> VisitedInternetLink --------------
$hyperlinks2 = $d.Hyperlinks
for ($i = 1; $i -le $hyperlinks2.Count; $i++) {
    $hl = $hyperlinks2.Item($i)
    if ($hl.Address -eq "https://fotos.eluniversal.com.mx/coleccion/muestra_fotogaleria.html?idgal=16498") {
        $hl.Range.Style = "VisitedInternetLink"
    }
}

Write-Output "done"
